$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on every touched cell first so that numeric-looking
# strings (e.g. "228.74", "38.815.99", "0.0838") are preserved exactly as
# text instead of being auto-coerced into floating point numbers.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '38.815.99'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.092.82'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +2.39%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '228.74'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.611'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '60.57'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.67%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.384'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +2.28%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0838'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.20%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.401.22'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +2.26%  '
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +4.51%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '21.90'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +2.12%  '
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +4.64%  '
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.091.56'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +2.13%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '38.702.78'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +2.59%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '71.65'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +3.30%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +2.20%  '
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +1.41%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '227.43'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +2.08%  '
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.39'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +3.39%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '171.16'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +1.31%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.48'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +1.82%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.140'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +9.03%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +13.22%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '19.16'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +2.08%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.98%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.38'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +5.50%  '
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +3.07%  '
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +4.75%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +1.72%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -1.61%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +1.60%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +3.46%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.18'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +0.95%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.541.90'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.66%  '
$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '100.91'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +3.28%  '
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0224'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +4.21%  '
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.80%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +3.13%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +8.16%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '4.11'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -1.56%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +2.91%  '
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.290.07'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +2.39%  '
